# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# The "Rules" sheet's last rule row (row 11) changes the value of the
# "Rule" name cell (B11) from the text "R40" to the text "1". The cell
# keeps its existing (string) type, so we write a text value rather than
# letting Excel auto-convert the digit into a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")
$cell.Value = "'1"
